$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-19 22:17:19"
$wsZhCn.Range("E5").Value = "2016-03-19 22:17:19"
$wsZhCn.Range("H4").Value = "2016-03-19 22:17:38"
$wsZhCn.Range("H5").Value = "2016-03-19 22:17:38"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-19 22:17:23"
$wsDeDe.Range("E5").Value = "2016-03-19 22:17:23"
$wsDeDe.Range("H4").Value = "2016-03-19 22:17:44"
$wsDeDe.Range("H5").Value = "2016-03-19 22:17:44"
